$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date text values (slash separators -> dash separators)
$dates = @{
    "A3"  = "28-07-2022"
    "A4"  = "01-08-2022"
    "A5"  = "04-08-2022"
    "A6"  = "08-08-2022"
    "A7"  = "11-08-2022"
    "A8"  = "15-08-2022"
    "A9"  = "18-08-2022"
    "A10" = "22-08-2022"
    "A11" = "25-08-2022"
    "A12" = "29-08-2022"
    "A13" = "01-09-2022"
    "A14" = "05-09-2022"
    "A15" = "08-09-2022"
    "A16" = "12-09-2022"
    "A17" = "15-09-2022"
    "A18" = "19-09-2022"
    "A19" = "22-09-2022"
    "A20" = "26-09-2022"
    "A21" = "29-09-2022"
}

foreach ($addr in $dates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$addr]
}

# Update the numeric attendance counters that changed
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("H14").Value = 0

$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("H15").Value = 0

$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("H20").Value = 0

$wb.Save()
